$d = $word.ActiveDocument

function Insert-LineBreak($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- "Programa" section, Portuguese paragraph ---
Insert-LineBreak "entes processos2. Processos bi" "entes processos^l2. Processos bi"
Insert-LineBreak "s de alimentos 3. Processos bi" "s de alimentos ^l3. Processos bi"
Insert-LineBreak "s desidratados 4. Discussão e " "s desidratados ^l4. Discussão e "
Insert-LineBreak "ncia industrial5. Bioenergia e" "ncia industrial^l5. Bioenergia e"

# --- "Programa" section, English paragraph ---
Insert-LineBreak "erent processes2. Biochemical " "erent processes^l2. Biochemical "
Insert-LineBreak "n/modifications3. Biochemical " "n/modifications^l3. Biochemical "
Insert-LineBreak "drated products4. Discussion a" "drated products^l4. Discussion a"
Insert-LineBreak "rial importance5. Bioenergy an" "rial importance^l5. Bioenergy an"

# --- "Bibliografia" section ---
Insert-LineBreak " 9788521313823.LIMA, U. A. Bio" " 9788521313823.^lLIMA, U. A. Bio"
Insert-LineBreak " 9788521214571.Moraes, I. O. B" " 9788521214571.^lMoraes, I. O. B"
